$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 669; this shifts the existing rows 669:710
# down to 670:711 (matches dimension change A1:D710 -> A1:D711).
$ws.Rows.Item(669).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds a date formatted as plain text (e.g. "2026/12/29") in this
# sheet, so force text entry to avoid Excel auto-converting the
# slash-delimited string into a date serial number.
$ws.Cells.Item(669, 1).NumberFormat = "@"
$ws.Cells.Item(669, 1).Value = "2026/01/21"
$ws.Cells.Item(669, 1).ClearFormats()

$ws.Cells.Item(669, 2).Value = "水"
$ws.Cells.Item(669, 3).Value = 14
$ws.Cells.Item(669, 4).Value = 25
